# Daily attendance processing - rotate "Recorded By" list entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val.Split(",") | ForEach-Object { $_.Trim() }
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
